$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$ws.Range("B2").Value = 0.40972222222222227
$ws.Range("B3").Value = 0.41666666666666669

$ws.Range("B3").Select()
